$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) following a data repull / recalculation.
$updates = @{
    2  = -5
    5  = 0
    7  = 2
    13 = -2
    15 = -3
    20 = -1
    22 = -3
    23 = -3
    25 = 4
    32 = -1
    33 = 2
    36 = 0
    37 = 1
    40 = 1
    51 = 1
    53 = -2
    56 = -2
    57 = -2
    59 = 5
    65 = 0
    69 = -4
    70 = -5
    72 = -2
    74 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
